# Hindalco price sheet update: a new day's row (13-02-2026) is published,
# pushing the whole history table down by one row. Everything that used to
# be row 2 now lives in row 3, row 3 -> row 4, etc. The brand new row 2
# re-uses the data that was previously in row 2 (same price/circular/link)
# but is now labeled with the new "as of" date 13-02-2026.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 2; this shifts the row 2..247
# cell values (and their number formats) down to rows 3..248.
$ws.Rows("2:2").Insert()

# The newly inserted row 2 is blank. Populate it by duplicating the row
# immediately below (which now holds what used to be the old row 2's
# data), preserving styles/number formats exactly.
$ws.Range("A3:F3").Copy($ws.Range("A2:F2"))

# Only the "Date" column actually changes for the new top row.
$ws.Range("A2").Value = "13-02-2026"

# The row-insert above shifts cell values, but (in this environment) it
# does not relocate the hyperlink objects attached to column F - they stay
# pinned to their original absolute rows. Rebuild the F-column hyperlinks
# from scratch so each row 2..248 points at the PDF link shown in its own
# cell text.
$ws.Cells.Hyperlinks.Delete()

$lastRow = 248
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
}
